# Repull data, push all data, mean calculation
# Updates the dSF (column F) values for specific rows to reflect the
# re-pulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = -4
    9  = -5
    12 = -6
    15 = -8
    16 = 4
    24 = 0
    25 = 2
    27 = 2
    28 = 0
    32 = -8
    37 = 0
    39 = 0
    46 = -7
    47 = -3
    50 = 1
    54 = -4
    58 = -10
    60 = -4
    62 = -5
    63 = -6
    66 = -2
    68 = 3
    70 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
